$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last refreshed" timestamp banner in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 04:22"

function Set-CountryRow($row, $name, $total, $new, $active, $recovered, $critical, $deathsToday, $deaths) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $total
    $ws.Cells.Item($row, 3).Value = $new
    $ws.Cells.Item($row, 4).Value = $active
    $ws.Cells.Item($row, 5).Value = $recovered
    $ws.Cells.Item($row, 6).Value = $critical
    $ws.Cells.Item($row, 7).Value = $deathsToday
    $ws.Cells.Item($row, 8).Value = $deaths
}

# --- Refresh the country table rows 98-103 (data re-sorted by total cases) ---
# Row 98: Bolivia jumps above Albania in the sort order with updated stats.
Set-CountryRow 98 "Bolivia" 598 34 37 527 3 1 34

# Row 99: now shows Albania, carrying Albania's previous (row 98) stats.
Set-CountryRow 99 "Albania" 584 0 327 231 5 0 26

# Row 100: now shows Burkina Faso, carrying its previous (row 99) stats.
Set-CountryRow 100 "Burkina Faso" 581 0 357 186 0 0 38

# Row 101: now shows Kirguistan, carrying its previous (row 100) stats.
Set-CountryRow 101 "Kirguistan" 568 0 201 360 5 0 7

# Row 102 (Uruguay) is unchanged.

# Row 103: Honduras stays in place but its stats are updated.
Set-CountryRow 103 "Honduras" 494 17 25 423 10 0 46
